# Apply updated crypto price/volume data to worksheet cells.
# Source cells are plain text (inline strings) so we force the
# NumberFormat to text ("@") before assigning values; this prevents
# Excel from auto-converting numeric-looking strings (e.g. "18.99")
# into floating point numbers and losing exact formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.284.76'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.595.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.99'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0854'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.819.77'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.595.05'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.98'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.503'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.44'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.275.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.06'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.91%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.92'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.30'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.41'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.37%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.483.27'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.62%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.569'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.76'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.99%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.17'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.932'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.732.62'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.755'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.14'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.65'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.09%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0956'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.03%  '
